$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # warp_execution_efficiency
$ws.Cells.Item(2, 11).Value = 74.842124
$ws.Cells.Item(3, 11).Value = 82.019994
$ws.Cells.Item(4, 11).Value = 82.043796
$ws.Cells.Item(5, 11).Value = 72.361273
$ws.Cells.Item(6, 11).Value = 90.62459800000001
$ws.Cells.Item(7, 11).Value = 78.518708
$ws.Cells.Item(8, 11).Value = 70.559577
$ws.Cells.Item(9, 11).Value = 81.149231
$ws.Cells.Item(10, 11).Value = 83.42728700000001
$ws.Cells.Item(11, 11).Value = 89.822239
$ws.Cells.Item(12, 11).Value = 89.772651
$ws.Cells.Item(13, 11).Value = 86.599085
$ws.Cells.Item(14, 11).Value = 86.755745
$ws.Cells.Item(15, 11).Value = 72.69129
$ws.Cells.Item(16, 11).Value = 81.07543200000001
$ws.Cells.Item(17, 11).Value = 87.79813300000001
$ws.Cells.Item(18, 11).Value = 87.000912
$ws.Cells.Item(19, 11).Value = 72.69544399999999
$ws.Cells.Item(20, 11).Value = 83.52371599999999
$ws.Cells.Item(21, 11).Value = 86.867188
$ws = $wb.Worksheets.Item(2)  # global_load_requests
$ws.Cells.Item(2, 11).Value = 1768835
$ws.Cells.Item(3, 11).Value = 363605
$ws.Cells.Item(4, 11).Value = 746579
$ws.Cells.Item(5, 11).Value = 344220
$ws.Cells.Item(6, 11).Value = 8156159
$ws.Cells.Item(7, 11).Value = 5750419
$ws.Cells.Item(8, 11).Value = 10285127
$ws.Cells.Item(9, 11).Value = 11919276
$ws.Cells.Item(10, 11).Value = 5202544
$ws.Cells.Item(11, 11).Value = 19783898
$ws.Cells.Item(12, 11).Value = 149115138
$ws.Cells.Item(13, 11).Value = 53444876
$ws.Cells.Item(14, 11).Value = 84606075
$ws.Cells.Item(16, 11).Value = 1866370052
$ws.Cells.Item(17, 11).Value = 487374612
$ws.Cells.Item(18, 11).Value = 539683532
$ws.Cells.Item(20, 11).Value = 16764941222
$ws.Cells.Item(21, 11).Value = 13066278928
$ws = $wb.Worksheets.Item(3)  # gld_transactions_per_request
$ws.Cells.Item(2, 11).Value = 3.069095
$ws.Cells.Item(3, 11).Value = 9.499979
$ws.Cells.Item(4, 11).Value = 9.964525
$ws.Cells.Item(5, 11).Value = 11.31592
$ws.Cells.Item(6, 11).Value = 4.380368
$ws.Cells.Item(7, 11).Value = 5.106219
$ws.Cells.Item(8, 11).Value = 3.572831
$ws.Cells.Item(9, 11).Value = 5.011258
$ws.Cells.Item(10, 11).Value = 11.667279
$ws.Cells.Item(11, 11).Value = 7.13977
$ws.Cells.Item(12, 11).Value = 3.549011
$ws.Cells.Item(13, 11).Value = 4.700175
$ws.Cells.Item(14, 11).Value = 4.183839
$ws.Cells.Item(15, 11).Value = 12.88468
$ws.Cells.Item(16, 11).Value = 3.117551
$ws.Cells.Item(17, 11).Value = 3.551654
$ws.Cells.Item(18, 11).Value = 3.838467
$ws.Cells.Item(19, 11).Value = 11.347727
$ws.Cells.Item(20, 11).Value = 3.550774
$ws.Cells.Item(21, 11).Value = 2.9905
$ws = $wb.Worksheets.Item(4)  # gld_efficiency
$ws.Cells.Item(2, 11).Value = 50.433933
$ws.Cells.Item(3, 11).Value = 32.780793
$ws.Cells.Item(4, 11).Value = 32.598064
$ws.Cells.Item(5, 11).Value = 28.504774
$ws.Cells.Item(6, 11).Value = 67.33591800000001
$ws.Cells.Item(7, 11).Value = 49.875182
$ws.Cells.Item(8, 11).Value = 42.792426
$ws.Cells.Item(9, 11).Value = 52.003516
$ws.Cells.Item(10, 11).Value = 31.194865
$ws.Cells.Item(11, 11).Value = 53.703274
$ws.Cells.Item(12, 11).Value = 70.250327
$ws.Cells.Item(13, 11).Value = 54.563376
$ws.Cells.Item(14, 11).Value = 56.021382
$ws.Cells.Item(15, 11).Value = 25.464429
$ws.Cells.Item(16, 11).Value = 73.978745
$ws.Cells.Item(17, 11).Value = 66.46365299999999
$ws.Cells.Item(18, 11).Value = 65.02746500000001
$ws.Cells.Item(19, 11).Value = 28.795671
$ws.Cells.Item(20, 11).Value = 77.564954
$ws.Cells.Item(21, 11).Value = 68.345524
